$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3971766666666667
$ws.Range("H2").Value = 1.19153
$ws.Range("Q2").Value = 0.08540741408555556
$ws.Range("R2").Value = 0.76866672677
